$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180 (shifts existing rows 180-244 down to 181-245,
# dimension grows from A1:R244 to A1:R245) and fill it with a new weekly
# price observation for Albahaca at Femacal de La Calera.
$ws.Rows.Item(180).Insert()

$ws.Range("A180").Value = 3
$ws.Range("B180").Value = "Femacal de La Calera"
$ws.Range("C180").Value = "Coquimbo"
$ws.Range("D180").Value = 44988
$ws.Range("E180").Value = 5
$ws.Range("F180").Value = 100112052
$ws.Range("G180").Value = "Albahaca"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 105
$ws.Range("K180").Value = 4500
$ws.Range("L180").Value = 5000
$ws.Range("M180").Value = 4738
$ws.Range("N180").Value = "`$/docena de matas"
$ws.Range("O180").Value = "Provincia de Quillota"
$ws.Range("P180").Value = 790
$ws.Range("Q180").Value = 6
$ws.Range("R180").Value = "Hortaliza"
